# edit.ps1 — apply the "Internship Final Report" proofing-pass edit.
#
# The authoritative diff splits several existing runs into multiple runs
# (wrapping the newly-isolated pieces with <w:proofErr .../> spell/grammar
# markers), removes a trailing empty paragraph, and makes two package-level
# tweaks (a VML autoshape id renumber and a <w:semiHidden/> flag on the
# DefaultParagraphFont style) that have no surface in the Word object model
# exposed here (no Shapes/InlineShapes for legacy VML <v:rect> "hr" picts,
# and Styles.Item(...).SemiHidden does not exist). Those two are therefore
# not reachable from COM automation and are skipped; every run-split and
# text-anchored change is applied below.
#
# Technique: toggling a character-formatting property (Bold) over a
# sub-range and then toggling it back is the only reliable way this COM
# surface exposes to force Word to materialize a dedicated <w:r> for that
# sub-range (Find/Replace alone keeps the whole match inside one run).

$d = $word.ActiveDocument

function Split-Run([int]$startPos, [int]$endPos) {
    $rng = $d.Range($startPos, $endPos)
    $rng.Font.Bold = 1
    $rng.Font.Bold = 0
}

# ---------------------------------------------------------------------
# 1) "Company: ShadowFox" -> "Company: " + [space run] + "ShadowFox"
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("ShadowFox")
Split-Run $r.Start $r.End

# ---------------------------------------------------------------------
# 2) Objectives paragraph:
#    "...HTML, CSS, JavaScript, " + " " + "and other web development..."
#    -> "...HTML, CSS, " + "JavaScript, " + " " + "and" + " other web development..."
# ---------------------------------------------------------------------
$r = $d.Content
$r.Find.Execute("JavaScript, ")
Split-Run $r.Start $r.End

$r = $d.Content
$r.Find.Execute("and other web development languages in real-world")
$r.Collapse(1)
$r.MoveEnd(1, 3)                      # just "and"
Split-Run $r.Start $r.End

# ---------------------------------------------------------------------
# 3) "During my internship at ShadowFox, my responsibilities included:"
#    -> "During my internship at " + "ShadowFox" + ", my responsibilities included:"
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(7).Range
$scope = $d.Range($p.Start, $p.End)
$scope.Find.Execute("ShadowFox")
Split-Run $scope.Start $scope.End

# ---------------------------------------------------------------------
# 4) Conclusion paragraph:
#    "This internship at ShadowFox provided ... professional workflows ."
#    "I am confident ... future career."
#    -> "This internship at " + "ShadowFox" + " provided ... professional "
#       + "workflows " + "." + "I" + " am confident ... future career."
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(23).Range
$scope = $d.Range($p.Start, $p.End)
$scope.Find.Execute("ShadowFox")
Split-Run $scope.Start $scope.End

$p = $d.Paragraphs.Item(23).Range
$scope = $d.Range($p.Start, $p.End)
$scope.Find.Execute("workflows ")
Split-Run $scope.Start $scope.End

$p = $d.Paragraphs.Item(23).Range
$scope = $d.Range($p.Start, $p.End)
$scope.Find.Execute("I am confident")
$scope.Collapse(1)
$scope.MoveEnd(1, 1)                  # just "I"
Split-Run $scope.Start $scope.End

# ---------------------------------------------------------------------
# 5) Acknowledgments paragraph:
#    "...entire ShadowFox team..." -> "...entire " + "ShadowFox" + " team..."
# ---------------------------------------------------------------------
$p = $d.Paragraphs.Item(25).Range
$scope = $d.Range($p.Start, $p.End)
$scope.Find.Execute("ShadowFox")
Split-Run $scope.Start $scope.End

# ---------------------------------------------------------------------
# 6) Drop the stray trailing empty paragraph (last paragraph in the body,
#    right before sectPr) by deleting its paragraph mark together with the
#    mark that ends the Acknowledgments paragraph.
# ---------------------------------------------------------------------
$last = $d.Paragraphs.Count
$ackPara = $d.Paragraphs.Item($last - 1)
$emptyPara = $d.Paragraphs.Item($last)
if ($emptyPara.Range.Text -eq "") {
    $delRange = $d.Range($ackPara.Range.End - 1, $emptyPara.Range.End)
    $delRange.Delete()
}

Write-Output "done"
